$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear old sample data (A1:A5) from the sheet
$ws.Range("A1:A5").Clear()

# ---- Header row (row 3) ----
$ws.Range("C3").Value = "Roll no"
$ws.Range("D3").Value = "Marks"
$ws.Range("E3").Value = "Pass  or fail"
$ws.Range("E3:F3").Merge()
$ws.Rows.Item(3).RowHeight = 18.75

# ---- Student rows (4-13): Roll no / Marks ----
$rolls = 101,102,103,104,105,106,107,108,109,110
$marks = 86,94,100,56,31,35,36,38,98,23

for ($i = 0; $i -lt 10; $i++) {
    $r = 4 + $i
    $ws.Range("C$r").Value = $rolls[$i]
    $ws.Range("D$r").Value = $marks[$i]
    $ws.Range("E$r").Formula = '=IF(D' + $r + '>35,"Pass","Fail")'
    $ws.Range("E$r`:F$r").Merge()
}

# ---- Total row (14) ----
$ws.Range("C14").Value = "Total"
$ws.Range("D14").Formula = "=SUM(D4:D13)"
$ws.Range("E14:F14").Merge()

$ws.Calculate()
